$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 20877.8
$ws.Range("J17").Value = 22184.895
$ws.Range("L17").Value = 66554.685
$ws.Range("N17").Value = -66890.685
$ws.Range("H34").Value = 2727.7856
$ws.Range("I34").Value = 2745.3845
$ws.Range("J34").Value = 2499
$ws.Range("K34").Value = 2745.3845
$ws.Range("L34").Value = 2499
$ws.Range("M34").Value = -2542.3845
$ws.Range("N34").Value = -2905
$ws.Range("H36").Value = 2727.7856
$ws.Range("I36").Value = 2745.3845
$ws.Range("J36").Value = 2499
$ws.Range("K36").Value = 2745.3845
$ws.Range("L36").Value = 2499
$ws.Range("M36").Value = -2030.3845
$ws.Range("N36").Value = -3929
$ws.Range("H40").Value = 100002984
$ws.Range("J40").Value = 100002984
$ws.Range("L40").Value = 100002984
$ws.Range("N40").Value = -100003334
$ws.Range("H74").Value = 5285.091
$ws.Range("I74").Value = 4329.625
$ws.Range("K74").Value = 4329.625
$ws.Range("M74").Value = -3393.625
$ws.Range("H76").Value = 9980.857
$ws.Range("I76").Value = 12535.8
$ws.Range("J76").Value = 3593.5
$ws.Range("K76").Value = 12535.8
$ws.Range("L76").Value = 3593.5
$ws.Range("M76").Value = -12220.8
$ws.Range("N76").Value = -4223.5
$ws.Range("H77").Value = 5285.091
$ws.Range("I77").Value = 4329.625
$ws.Range("K77").Value = 21648.125
$ws.Range("M77").Value = -16968.125
$ws.Range("H79").Value = 9980.857
$ws.Range("I79").Value = 12535.8
$ws.Range("J79").Value = 3593.5
$ws.Range("K79").Value = 12535.8
$ws.Range("L79").Value = 3593.5
$ws.Range("M79").Value = -11443.8
$ws.Range("N79").Value = -5777.5
$ws.Range("H100").Value = 12032.833
$ws.Range("I100").Value = 12375
$ws.Range("J100").Value = 11861.75
$ws.Range("K100").Value = 12375
$ws.Range("L100").Value = 11861.75
$ws.Range("M100").Value = -11834
$ws.Range("N100").Value = -12943.75
$ws.Range("H107").Value = 1795.2667
$ws.Range("I107").Value = 524.3
$ws.Range("K107").Value = 524.3
$ws.Range("M107").Value = 1395.7
$ws.Range("H137").Value = 1574.6072
$ws.Range("I137").Value = 995.3684
$ws.Range("J137").Value = 2797.4443
$ws.Range("K137").Value = 2986.1052
$ws.Range("L137").Value = 8392.332900000001
$ws.Range("M137").Value = -436.1052
$ws.Range("N137").Value = -13492.3329
$ws.Range("H138").Value = 2610.5132
$ws.Range("I138").Value = 1394
$ws.Range("J138").Value = 4582.1035
$ws.Range("K138").Value = 4182
$ws.Range("L138").Value = 13746.3105
$ws.Range("M138").Value = 958
$ws.Range("N138").Value = -24026.3105

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1376.64
$ws.Range("I97").Value = 1462.2273
$ws.Range("J97").Value = 749
$ws.Range("K97").Value = 1462.2273
$ws.Range("L97").Value = 749
$ws.Range("M97").Value = -966.2273
$ws.Range("N97").Value = -1741
$ws.Range("H102").Value = 4453.7144
$ws.Range("I102").Value = 4033.0833
$ws.Range("K102").Value = 4033.0833
$ws.Range("M102").Value = -2411.0833
$ws.Range("H132").Value = 1926462.4
$ws.Range("I132").Value = 3328.5745
$ws.Range("K132").Value = 9985.7235
$ws.Range("M132").Value = -7455.7235

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2926.6365
$ws.Range("I94").Value = 2905.0588
$ws.Range("J94").Value = 3000
$ws.Range("K94").Value = 2905.0588
$ws.Range("L94").Value = 3000
$ws.Range("M94").Value = -2454.0588
$ws.Range("N94").Value = -3902
$ws.Range("H105").Value = 647549
$ws.Range("I105").Value = 992421.7
$ws.Range("K105").Value = 992421.7
$ws.Range("M105").Value = -990674.7
$ws.Range("H134").Value = 2943299.2
$ws.Range("I134").Value = 2051.1724
$ws.Range("K134").Value = 6153.5172
$ws.Range("M134").Value = -3618.5172

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4025.5454
$ws.Range("I134").Value = 4142.5557
$ws.Range("J134").Value = 3499
$ws.Range("K134").Value = 12427.6671
$ws.Range("L134").Value = 10497
$ws.Range("M134").Value = -9892.667099999999
$ws.Range("N134").Value = -15567

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3750
$ws.Range("I68").Value = 3000
$ws.Range("J68").Value = 4000
$ws.Range("K68").Value = 9000
$ws.Range("L68").Value = 12000
$ws.Range("M68").Value = -8189
$ws.Range("N68").Value = -13622
$ws.Range("H71").Value = 3750
$ws.Range("I71").Value = 3000
$ws.Range("J71").Value = 4000
$ws.Range("K71").Value = 27000
$ws.Range("L71").Value = 36000
$ws.Range("M71").Value = -22944
$ws.Range("N71").Value = -44112
$ws.Range("H121").Value = 5087.375
$ws.Range("I121").Value = 748.5
$ws.Range("J121").Value = 5707.2144
$ws.Range("K121").Value = 2245.5
$ws.Range("L121").Value = 17121.6432
$ws.Range("M121").Value = -935.5
$ws.Range("N121").Value = -19741.6432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1388
$ws.Range("I97").Value = 1624.75
$ws.Range("J97").Value = 1198.6
$ws.Range("K97").Value = 1624.75
$ws.Range("L97").Value = 1198.6
$ws.Range("M97").Value = -1128.75
$ws.Range("N97").Value = -2190.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3708591
$ws.Range("I93").Value = 1878.7778
$ws.Range("J93").Value = 9268659
$ws.Range("K93").Value = 1878.7778
$ws.Range("L93").Value = 9268659
$ws.Range("M93").Value = -630.7778000000001
$ws.Range("N93").Value = -9271155
$ws.Range("H100").Value = 35754020
$ws.Range("J100").Value = 50055388
$ws.Range("L100").Value = 50055388
$ws.Range("N100").Value = -50056470
$ws.Range("H122").Value = 3392.8728
$ws.Range("I122").Value = 3396.255
$ws.Range("J122").Value = 3349.75
$ws.Range("K122").Value = 10188.765
$ws.Range("L122").Value = 10049.25
$ws.Range("M122").Value = -7738.764999999999
$ws.Range("N122").Value = -14949.25
$ws.Range("H136").Value = 2397.8438
$ws.Range("I136").Value = 1834.16
$ws.Range("K136").Value = 5502.48
$ws.Range("M136").Value = -2952.48

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 7798
$ws.Range("I96").Value = 10150.25
$ws.Range("J96").Value = 6229.8335
$ws.Range("K96").Value = 10150.25
$ws.Range("L96").Value = 6229.8335
$ws.Range("M96").Value = -8777.25
$ws.Range("N96").Value = -8975.833500000001
$ws.Range("H122").Value = 3297.875
$ws.Range("I122").Value = 2799.2856
$ws.Range("J122").Value = 6788
$ws.Range("K122").Value = 8397.856800000001
$ws.Range("L122").Value = 20364
$ws.Range("M122").Value = -5947.856800000001
$ws.Range("N122").Value = -25264
$ws.Range("H132").Value = 417855.22
$ws.Range("I132").Value = 1160.2273
$ws.Range("J132").Value = 5001500
$ws.Range("K132").Value = 3480.6819
$ws.Range("L132").Value = 15004500
$ws.Range("M132").Value = -950.6819
$ws.Range("N132").Value = -15009560
